# Auto-generated Excel COM-interop script
# Applies updated currentAveragePrice / LevePrice / LeveProfit values
# sourced from the scheduled market-data refresh across all job sheets.

$wb = $excel.ActiveWorkbook

# --- ALC sheet updates ---
$ws = $wb.Worksheets.Item("ALC")
$ALCUpdates = @(
    @("H100", 1575),
    @("I100", 1500),
    @("K100", 1500),
    @("M100", -959),
    @("H115", 9091544),
    @("I115", 9091544),
    @("K115", 27274632),
    @("M115", -27273065),
    @("H132", 2777.7144),
    @("I132", 2825.2727),
    @("J132", 1993),
    @("K132", 8475.8181),
    @("L132", 5979),
    @("M132", -5945.8181),
    @("N132", -11039),
    @("H137", 35860200),
    @("I137", 50002880),
    @("J137", 503501.5),
    @("K137", 150008640),
    @("L137", 1510504.5),
    @("M137", -150006090),
    @("N137", -1515604.5)
)
foreach ($u in $ALCUpdates) {
    $ws.Range($u[0]).Value = $u[1]
}

# --- ARM sheet updates ---
$ws = $wb.Worksheets.Item("ARM")
$ARMUpdates = @(
    @("H32", 12503416),
    @("I32", 14085989),
    @("J32", 18668.223),
    @("K32", 14085989),
    @("L32", 18668.223),
    @("M32", -14085702),
    @("N32", -19242.223),
    @("H61", 2093.258),
    @("I61", 2248.7727),
    @("J61", 1713.1111),
    @("K61", 2248.7727),
    @("L61", 1713.1111),
    @("M61", -2036.7727),
    @("N61", -2137.1111),
    @("H102", 3680.2),
    @("I102", 2095),
    @("K102", 2095),
    @("M102", -473),
    @("H132", 1987.5385),
    @("I132", 2048.077),
    @("J132", 1927),
    @("K132", 6144.231000000001),
    @("L132", 5781),
    @("M132", -3614.231000000001),
    @("N132", -10841),
    @("H136", 2093.258),
    @("I136", 2248.7727),
    @("J136", 1713.1111),
    @("K136", 6746.3181),
    @("L136", 5139.3333),
    @("M136", -4196.3181),
    @("N136", -10239.3333)
)
foreach ($u in $ARMUpdates) {
    $ws.Range($u[0]).Value = $u[1]
}

# --- BSM sheet updates ---
$ws = $wb.Worksheets.Item("BSM")
$BSMUpdates = @(
    @("H107", 1271.2),
    @("I107", 1022.3158),
    @("K107", 1022.3158),
    @("M107", 897.6842),
    @("H134", 1298.6875),
    @("I134", 1369.76),
    @("J134", 1044.8572),
    @("K134", 4109.28),
    @("L134", 3134.5716),
    @("M134", -1574.28),
    @("N134", -8204.571599999999)
)
foreach ($u in $BSMUpdates) {
    $ws.Range($u[0]).Value = $u[1]
}

# --- CRP sheet updates ---
$ws = $wb.Worksheets.Item("CRP")
$CRPUpdates = @(
    @("H31", 1350.3778),
    @("I31", 1817.5714),
    @("J31", 1139.3871),
    @("K31", 1817.5714),
    @("L31", 1139.3871),
    @("M31", -1522.5714),
    @("N31", -1729.3871),
    @("H34", 1350.3778),
    @("I34", 1817.5714),
    @("J34", 1139.3871),
    @("K34", 1817.5714),
    @("L34", 1139.3871),
    @("M34", -1615.5714),
    @("N34", -1543.3871),
    @("H132", 1556.1786),
    @("I132", 1193.7),
    @("J132", 2462.375),
    @("K132", 3581.1),
    @("L132", 7387.125),
    @("M132", -1051.1),
    @("N132", -12447.125),
    @("H134", 3615.0952),
    @("I134", 1248.0769),
    @("J134", 7461.5),
    @("K134", 3744.2307),
    @("L134", 22384.5),
    @("M134", -1209.2307),
    @("N134", -27454.5)
)
foreach ($u in $CRPUpdates) {
    $ws.Range($u[0]).Value = $u[1]
}

# --- CUL sheet updates ---
$ws = $wb.Worksheets.Item("CUL")
$CULUpdates = @(
    @("H6", 92.125),
    @("I6", 92.125),
    @("K6", 276.375),
    @("M6", -163.375),
    @("H18", 570),
    @("I18", 303),
    @("J18", 970.5),
    @("K18", 909),
    @("L18", 2911.5),
    @("M18", -740),
    @("N18", -3249.5),
    @("H68", 1341.4186),
    @("I68", 1409.3334),
    @("J68", 1305.0358),
    @("K68", 4228.0002),
    @("L68", 3915.1074),
    @("M68", -3417.0002),
    @("N68", -5537.107400000001),
    @("H71", 1341.4186),
    @("I71", 1409.3334),
    @("J71", 1305.0358),
    @("K71", 12684.0006),
    @("L71", 11745.3222),
    @("M71", -8628.000599999999),
    @("N71", -19857.3222),
    @("H100", 3645),
    @("J100", 4000),
    @("L100", 12000),
    @("N100", -13622),
    @("H113", 1078009),
    @("I113", 2155540.8),
    @("J113", 477.125),
    @("K113", 6466622.399999999),
    @("L113", 1431.375),
    @("M113", -6464452.399999999),
    @("N113", -5771.375),
    @("H131", 828.3),
    @("J131", 1033.5862),
    @("L131", 3100.7586),
    @("N131", -13180.7586)
)
foreach ($u in $CULUpdates) {
    $ws.Range($u[0]).Value = $u[1]
}

# --- GSM sheet updates ---
$ws = $wb.Worksheets.Item("GSM")
$GSMUpdates = @(
    @("H132", 3958.2334),
    @("I132", 3780.0417),
    @("J132", 4671),
    @("K132", 11340.1251),
    @("L132", 14013),
    @("M132", -8810.125100000001),
    @("N132", -19073)
)
foreach ($u in $GSMUpdates) {
    $ws.Range($u[0]).Value = $u[1]
}

# --- LTW sheet updates ---
$ws = $wb.Worksheets.Item("LTW")
$LTWUpdates = @(
    @("H22", 45909988),
    @("I22", 63125844),
    @("J22", 1043.1666),
    @("K22", 63125844),
    @("L22", 1043.1666),
    @("M22", -63125549),
    @("N22", -1633.1666),
    @("H27", 45909988),
    @("I27", 63125844),
    @("J27", 1043.1666),
    @("K27", 63125844),
    @("L27", 1043.1666),
    @("M27", -63125737),
    @("N27", -1257.1666),
    @("H136", 1626.2),
    @("I136", 1616.3478),
    @("J136", 1658.5714),
    @("K136", 4849.0434),
    @("L136", 4975.7142),
    @("M136", -2299.0434),
    @("N136", -10075.7142),
    @("H141", 42715),
    @("J141", 42715),
    @("L141", 42715),
    @("N141", -53075)
)
foreach ($u in $LTWUpdates) {
    $ws.Range($u[0]).Value = $u[1]
}

# --- WVR sheet updates ---
$ws = $wb.Worksheets.Item("WVR")
$WVRUpdates = @(
    @("H132", 1601.7576),
    @("I132", 1371.2354),
    @("J132", 1846.6875),
    @("K132", 4113.706200000001),
    @("L132", 5540.0625),
    @("M132", -1583.706200000001),
    @("N132", -10600.0625)
)
foreach ($u in $WVRUpdates) {
    $ws.Range($u[0]).Value = $u[1]
}
